# Update transition-probability matrix values on Sheet1 of the active workbook.
# These cells represent probabilities recomputed after adding more games /
# speeding up the simulate-game logic (see commit message). Each affected
# row still sums to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("B2").Value = 0.125
$ws.Range("C2").Value = 0.625
$ws.Range("P2").Value = 0.125
$ws.Range("S2").Value = 0.125

$ws.Range("C3").Value = 0.2857142857142857
$ws.Range("P3").Value = 0.4285714285714285
$ws.Range("S3").Value = 0.2857142857142857

$ws.Range("P4").Value = 0.8333333333333334
$ws.Range("S4").Value = 0.1666666666666667

$ws.Range("B6").Value = 0.1
$ws.Range("J6").Value = 0.2
$ws.Range("S6").Value = 0.7

$ws.Range("J7").Value = 0.1428571428571428
$ws.Range("Q7").Value = 0.2857142857142857
$ws.Range("R7").Value = 0.2857142857142857
$ws.Range("S7").Value = 0.2857142857142857

$ws.Range("F8").Value = 0.1666666666666667
$ws.Range("J8").Value = 0.2777777777777778
$ws.Range("O8").Value = 0.05555555555555555
$ws.Range("Q8").Value = 0.1111111111111111
$ws.Range("R8").Value = 0.2222222222222222
$ws.Range("S8").Value = 0.1666666666666667

$ws.Range("F9").Value = 0.09090909090909091
$ws.Range("J9").Value = 0.2727272727272727
$ws.Range("Q9").Value = 0.09090909090909091
$ws.Range("R9").Value = 0.09090909090909091
$ws.Range("S9").Value = 0.4545454545454545

$ws.Range("B10").Value = 0.04347826086956522
$ws.Range("D10").Value = 0.04347826086956522
$ws.Range("F10").Value = 0.02173913043478261
$ws.Range("J10").Value = 0.1521739130434783
$ws.Range("O10").Value = 0.01449275362318841
$ws.Range("Q10").Value = 0.3405797101449275
$ws.Range("R10").Value = 0.08695652173913043
$ws.Range("S10").Value = 0.2971014492753623

$ws.Range("G11").Value = 0.07142857142857142
$ws.Range("J11").Value = 0.1428571428571428
$ws.Range("K11").Value = 0.1428571428571428
$ws.Range("L11").Value = 0.6428571428571429

$ws.Range("G12").Value = 0.6
$ws.Range("J12").Value = 0.2
$ws.Range("L12").Value = 0.1
$ws.Range("S12").Value = 0.1

$ws.Range("H15").Value = 0.2941176470588235
$ws.Range("I15").Value = 0.05882352941176471
$ws.Range("J15").Value = 0.2352941176470588
$ws.Range("K15").Value = 0.05882352941176471
$ws.Range("O15").Value = 0.05882352941176471
$ws.Range("S15").Value = 0.2941176470588235

$ws.Range("J16").Value = 0.75
$ws.Range("O16").Value = 0.125
$ws.Range("S16").Value = 0.125

$ws.Range("F17").Value = 0.01923076923076923
$ws.Range("H17").Value = 0.1153846153846154
$ws.Range("I17").Value = 0.03846153846153846
$ws.Range("J17").Value = 0.5961538461538461
$ws.Range("K17").Value = 0.01923076923076923
$ws.Range("O17").Value = 0.0576923076923077
$ws.Range("S17").Value = 0.1538461538461539

$ws.Range("H18").Value = 0.1
$ws.Range("I18").Value = 0.1
$ws.Range("J18").Value = 0.5
$ws.Range("K18").Value = 0.05
$ws.Range("O18").Value = 0.15
$ws.Range("S18").Value = 0.1

$ws.Range("H19").Value = 0.06666666666666667
$ws.Range("I19").Value = 0.06666666666666667
$ws.Range("J19").Value = 0.5777777777777777
$ws.Range("K19").Value = 0.1
$ws.Range("O19").Value = 0.04444444444444445
$ws.Range("S19").Value = 0.1444444444444444
